$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("D2").Value = -0.209
$ws.Range("G2").Value = -0.2652925531914894
$ws.Range("H2").Value = -0.2652925531914894
$ws.Range("I2").Value = -0.1582446808510638
$ws.Range("J2").Value = -0.1582446808510638
$ws.Range("K2").Value = -65.2
$ws.Range("L2").Value = -0.4335106382978723
$ws.Range("U2").Value = 8.279999999999999
$ws.Range("V2").Value = 0.7596330275229357
$ws.Range("W2").Value = 21.80602006688963
$ws.Range("X2").Value = 0.4264119303569793
$ws.Range("Y2").Value = 21.37960813653265
$ws.Range("Z2").Value = 2.557388199285836
$ws.Range("AA2").Value = -0.4046930794082639
$ws.Range("AB2").Value = 0.08533830539311128
$ws.Range("AC2").Value = -0.4900313848013751
$ws.Range("AD2").Value = 73.2
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 73.2
$ws.Range("AG2").Value = 64.92
$ws.Range("AH2").Value = 0.8703923900118905
$ws.Range("AI2").Value = 1.224080267558528
$ws.Range("AJ2").Value = 0.8562384595093642
$ws.Range("AK2").Value = 1.260093167701863
$ws.Range("AL2").Value = 21.4
$ws.Range("AM2").Value = 19.14
$ws.Range("AN2").Value = -2.64259927797834
$ws.Range("AO2").Value = -1.11214953271028
$ws.Range("AP2").Value = -2.343682310469314
$ws.Range("AQ2").Value = -1.243469174503657

# --- Row 3 updates ---
$ws.Range("B3").Value = "ATMA Participações S.A. (BOVESPA:ATMP3)"
$ws.Range("D3").Value = -0.209
$ws.Range("G3").Value = -0.2652925531914894
$ws.Range("H3").Value = -0.2652925531914894
$ws.Range("I3").Value = -0.1582446808510638
$ws.Range("J3").Value = -0.1582446808510638
$ws.Range("K3").Value = -65.2
$ws.Range("L3").Value = -0.4335106382978723
$ws.Range("U3").Value = 8.279999999999999
$ws.Range("V3").Value = 0.7596330275229357
$ws.Range("W3").Value = 21.80602006688963
$ws.Range("X3").Value = 0.4264119303569793
$ws.Range("Y3").Value = 21.37960813653265
$ws.Range("Z3").Value = 2.557388199285836
$ws.Range("AA3").Value = -0.4046930794082639
$ws.Range("AB3").Value = 0.08533830539311128
$ws.Range("AC3").Value = -0.4900313848013751
$ws.Range("AD3").Value = 73.2
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 73.2
$ws.Range("AG3").Value = 64.92
$ws.Range("AH3").Value = 0.8703923900118905
$ws.Range("AI3").Value = 1.224080267558528
$ws.Range("AJ3").Value = 0.8562384595093642
$ws.Range("AK3").Value = 1.260093167701863
$ws.Range("AL3").Value = 21.4
$ws.Range("AM3").Value = 19.14
$ws.Range("AN3").Value = -2.64259927797834
$ws.Range("AO3").Value = -1.11214953271028
$ws.Range("AP3").Value = -2.343682310469314
$ws.Range("AQ3").Value = -1.243469174503657
